$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add sample-code submission columns: Batsman / Bowlers / the script note
$ws.Range("G1").Value = "Batsman"
$ws.Range("H1").Value = "Bowlers"
$ws.Range("I1").Value = "Get list of bowlers from script"

# Leave the active cell on the newly added header cell, matching the saved selection
$ws.Range("I1").Select() | Out-Null
